$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 3
$ws.Range("E3").Value = 5
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Update selection to D3 (reflects <selection activeCell="D3" sqref="D3"/> in sheetView)
$ws.Range("D3").Select()
